$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Add missing values to existing row 6
$ws.Range("X6").Value = 0.29000100000000373
$ws.Range("Y6").Value = "Up"

# Add new row 7 data (repeater scan result)
$ws.Range("A7").Value = 42648.890590277777
$ws.Range("B7").Value = -4
$ws.Range("C7").Value = "Neutral"
$ws.Range("D7").Value = 42
$ws.Range("E7").Value = 16969
$ws.Range("F7").Value = 1057
$ws.Range("G7").Value = 58
$ws.Range("H7").Value = 39
$ws.Range("I7").Value = 86
$ws.Range("J7").Value = 13
$ws.Range("K7").Value = 40305
$ws.Range("L7").Value = 168
$ws.Range("M7").Value = 114
$ws.Range("N7").Value = 78
$ws.Range("O7").Value = 12
$ws.Range("P7").Value = "Named"
$ws.Range("Q7").Value = 38.916275631518758
$ws.Range("R7").Value = 0
$ws.Range("S7").Value = -0.0078
$ws.Range("S7").NumberFormat = "0.00%"
$ws.Range("T7").Value = -0.0305
$ws.Range("T7").NumberFormat = "0.00%"
$ws.Range("U7").Value = 14.62
$ws.Range("V7").Value = "N/A"
$ws.Range("W7").Value = -2
